$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.3586222133355363, 0.05906290126593206, 0.1328761516594792, 3.198037250614647, 2.470281004012463, 2.005870798417163, 0.2391051643663751, 2.403337428066493, 0.1061377151953922, 0.7142530321728202, 2.16557970365297)
    3 = @(0.3567354971009422, 0.05841900120496035, 0.1331146532641405, 3.202799337409928, 2.472147669621577, 2.013821287754226, 0.2403925642453224, 2.293932862647125, 0.1064719773979412, 0.6948529172636455, 2.189582963661733)
    4 = @(0.3557323763221945, 0.05803562181423771, 0.1333013751575702, 3.207321764178047, 2.474676461168073, 2.019602237705172, 0.2412636879096937, 2.227809643418453, 0.106691783417495, 0.6832563537928635, 2.205049848110636)
    5 = @(0.3553626942099584, 0.05788242739141225, 0.1333876143550778, 3.209566357754198, 2.476054137689005, 2.022184089318699, 0.2416389738834717, 2.201129000290138, 0.106785029259461, 0.6786100943881834, 2.211535968366158)
    6 = @(0.3553036713423978, 0.05785717363396614, 0.133402547831091, 3.209963322415675, 2.476303853602275, 2.022626455690016, 0.2417025161003927, 2.196714731601048, 0.1068007347994557, 0.6778433895372018, 2.212624050416384)
    7 = @(0.3557272323028542, 0.05803354345685108, 0.1333024970923855, 3.207350409601275, 2.474693635980813, 2.019636142226403, 0.2412686669503614, 2.227448744644136, 0.1066930260781218, 0.6831933708492954, 2.20513658030135)
    8 = @(0.3579394562028284, 0.05883841250115651, 0.1329500382234219, 3.199347272398938, 2.470637348977263, 2.00842542633508, 0.2395323285694211, 2.365396824423669, 0.1062499512506818, 0.7074985541813561, 2.173704757732803)
    9 = @(0.3635090599126158, 0.06051086411927287, 0.1325776941186483, 3.196352801242128, 2.473679664158226, 1.993581634078851, 0.2367667435627823, 2.644244503605535, 0.1054962283342924, 0.75765730650172, 2.117849670806686)
    10 = @(0.3683513434346111, 0.06179589966745169, 0.1324975529445034, 3.201922655869055, 2.482659249640591, 1.987037441788772, 0.2351239295648213, 2.854201909799542, 0.1050120605469562, 0.7960297043641873, 2.080336993735539)
    11 = @(0.3707171350468172, 0.06239250187638135, 0.1325029046140322, 3.206150174743271, 2.488217933320499, 1.985009464268245, 0.2344609142977454, 2.950825500876306, 0.1048067842599858, 0.8138166690669522, 2.06403653819245)
    12 = @(0.3716364219886827, 0.0626201289845838, 0.1325109261511415, 3.207995047989158, 2.490535477881764, 1.984378140212698, 0.23422196054549, 2.987574167970536, 0.1047311948011975, 0.8205996647451599, 2.0579739266914)
    13 = @(0.3714373962336879, 0.06257102987574825, 0.1325089322098343, 3.207586862254928, 2.490026887560731, 1.984508028003006, 0.2342728847804487, 2.979652604946693, 0.1047473791219122, 0.819136717107142, 2.059274723530942)
    14 = @(0.3707922962617545, 0.06241119482339741, 0.1325034444990898, 3.20629706081435, 2.488404334657474, 1.984954785710897, 0.2344410126547416, 2.953845644536386, 0.1048005225475617, 0.8143737599933161, 2.0635355578793)
    15 = @(0.3704002021530926, 0.0623135128202037, 0.1325008633396614, 3.205538809944713, 2.487438179520325, 1.985246235435312, 0.2345455734349358, 2.938058877779326, 0.1048333534038122, 0.8114624869055831, 2.066159773093982)
    16 = @(0.3682000112006136, 0.06175715042651575, 0.1324980430314078, 3.201680498038755, 2.482325689727844, 1.987189084989836, 0.2351689550092786, 2.847909671672653, 0.1050257763779845, 0.7948739342829541, 2.0814176555124)
    17 = @(0.3668919964926545, 0.06141890524553872, 0.1325070087987541, 3.199747674110426, 2.47956728992898, 1.988624136828918, 0.2355729673028399, 2.792890646786248, 0.1051476504037225, 0.7847820971796295, 2.090973697265863)
    18 = @(0.3661550061564611, 0.06122548956656715, 0.1325161018689194, 3.198795363031337, 2.478119414947599, 1.989538855950116, 0.2358132804988351, 2.76134994786662, 0.1052191590856868, 0.7790087200673241, 2.096542006961759)
    19 = @(0.3659081103007935, 0.06116019782903948, 0.1325198571935289, 3.198500290294604, 2.477652988491656, 1.98986389928146, 0.2358960096788287, 2.750688826118619, 0.1052436131363521, 0.7770593141307174, 2.098439690993695)
    20 = @(0.3670296489361533, 0.06145479482893279, 0.1325056471202473, 3.199936926285574, 2.479846568368259, 1.988462128203849, 0.2355291382168332, 2.79873667524987, 0.1051345308449516, 0.7858531645875573, 2.089948995300714)
    21 = @(0.370981142557099, 0.0624580961174459, 0.1325048938052547, 3.206669281003059, 2.48887514281725, 1.984819852817139, 0.2343913006648783, 2.961421452926857, 0.1047848549137926, 0.8157714692349742, 2.062281060725786)
    22 = @(0.3737001286987436, 0.06312374928985776, 0.1325393370469143, 3.212491715572682, 2.4960152963219, 1.98323586039524, 0.2337182745308937, 3.068674246796206, 0.1045688154723807, 0.8356013603890347, 2.044839791050796)
    23 = @(0.3722364778180349, 0.06276757606720906, 0.1325177629931176, 3.20925386970967, 2.492090832120994, 1.98400833774798, 0.2340710223009062, 3.011346620895949, 0.1046829795526083, 0.8249925229466726, 2.054089802470958)
    24 = @(0.3669673695288083, 0.06143856588896313, 0.1325062504660419, 3.199850870465795, 2.479719876885326, 1.9885350929396, 0.2355489282986198, 2.796093403928239, 0.1051404577062041, 0.7853688460569046, 2.090412031271303)
    25 = @(0.3618705479088788, 0.0600484600235589, 0.1326443890149687, 3.195800454670547, 2.471675208029865, 1.996831978285343, 0.2374465341802541, 2.567916001242793, 0.1056878658863578, 0.7438209092506725, 2.132341031889418)
}

$cols = @("C","D","E","F","G","H","J","K","L","M","N")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
